$d = $word.ActiveDocument

# The target sentence is unique in the document, so anchor on it directly.
$old = "Comunicația va fi non-blocantă, un mesaj va fi așteptat maxim 5 secunde"
$new = "Comunicația va fi non-blocantă, un mesaj va fi așteptat maxim 15 secunde"

$anchor = $d.Content
$found = $anchor.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) {
    throw "Could not find the target sentence to update."
}

# Re-locate "...maxim 1" to know exactly where the new "1" character ends, so we
# can split the run there (Word merges same-formatted runs on a text edit, so we
# restore the 3-way run split the source diff expects: "...maxim ", "1", "5 secunde").
$prefix = $d.Content
$prefix.Find.Execute("Comunicația va fi non-blocantă, un mesaj va fi așteptat maxim 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$posAfterOne = $prefix.End

# "1" becomes its own run: toggle a formatting property on/off to force Word to
# break it out of the surrounding run without altering its visible formatting.
$oneRange = $d.Range($posAfterOne - 1, $posAfterOne)
$oneRange.Font.Bold = 1
$oneRange.Font.Bold = 0

# "5 secunde" becomes its own run the same way; its end lines up exactly with the
# following (untouched) run so that run is left alone.
$fiveSecRange = $d.Range($posAfterOne, $posAfterOne + 9)
$fiveSecRange.Font.Bold = 1
$fiveSecRange.Font.Bold = 0
